$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds item ids as text (e.g. "1", "500"); force Text format
# before assigning so Excel does not silently convert them to numbers.
$ws.Range("A3:A14").NumberFormat = "@"

# Row 3
$ws.Range("A3").Value = '1'
$ws.Range("B3").Value = 'Alpha 21.1'
$ws.Range("C3").Value = '古旧的护身符'
$ws.Range("D3").Value = 'Old Charm'
$ws.Range("E3").Value = '古びたお守り'
$ws.Range("F3").Value = '有些似曾相识的褪色贝壳护身符。'
$ws.Range("G3").Value = 'A somewhat familiar looking amulet made of faded shells.'
$ws.Range("H3").Value = 'どこか見覚えのある、色あせた貝殻のお守りだ。'

# Row 4
$ws.Range("A4").Value = '2'
$ws.Range("B4").Value = 'Alpha 21.1'
$ws.Range("C4").Value = '用惯了的背包'
$ws.Range("D4").Value = 'Well-worn Backpack'
$ws.Range("E4").Value = '使い慣れたバックパック'
$ws.Range("F4").Value = '无论什么时候都很可靠的旅行伙伴。'
$ws.Range("G4").Value = 'Always reliable, no matter the journey.'
$ws.Range("H4").Value = 'どんな時でも頼りになる旅の相棒だ。'

# Row 5
$ws.Range("A5").Value = '3'
$ws.Range("B5").Value = 'Alpha 21.1'
$ws.Range("C5").Value = '女神的余香'
$ws.Range("D5").Value = 'Fragrance of Goddess '
$ws.Range("E5").Value = '女神の残り香'
$ws.Range("F5").Value = '幸运女神的轻柔余香。据说会吸引愿望女神。'
$ws.Range("G5").Value = 'The soft remaining fragrance of the Goddess of Luck. It is believed to draw the Goddess of Wish. '
$ws.Range("H5").Value = '幸運の女神のふわふわな残り香。願いの女神を引き寄せるといわれている。'

# Row 6
$ws.Range("A6").Value = '4'
$ws.Range("B6").Value = 'Beta 22.57'
$ws.Range("C6").Value = '黑猫的口水'
$ws.Range("D6").Value = 'Black Cat''s Saliva'
$ws.Range("E6").Value = '黒猫のよだれ'
$ws.Range("F6").Value = '含有强烈祝福魔力的唾液。

[持有效果]
加强水井许愿的效果。'
$ws.Range("G6").Value = 'Saliva containing the magical power of a strong blessing.

[Passive Effect]
Enhances the power of wishes made at wells.'
$ws.Range("H6").Value = '強い祝福の魔力が含まれた唾液だ。

[所持効果]
井戸の願いを強化する。'

# Row 7
$ws.Range("A7").Value = '5'
$ws.Range("B7").Value = 'EA 23.8'
$ws.Range("C7").Value = '女神的碎羽'
$ws.Range("D7").Value = 'Feather Shard of Goddess'
$ws.Range("E7").Value = '女神の羽くず'
$ws.Range("F7").Value = '治愈女神的翅膀上掉下的一片羽毛'
$ws.Range("G7").Value = 'A piece of feather fallen from the wing of the Goddess of Healing.'
$ws.Range("H7").Value = '癒しの女神の翼から抜け落ちた羽の一辺。'

# Row 8
$ws.Range("A8").Value = '500'
$ws.Range("B8").Value = 'Beta 22.57'
$ws.Range("C8").Value = '韦尔尼斯矿场的钥匙'
$ws.Range("D8").Value = 'Key for Vernis Mine'
$ws.Range("E8").Value = 'ヴェルニース炭鉱の鍵'
$ws.Range("F8").Value = '通往韦尔尼斯矿场深层之门的钥匙。'
$ws.Range("G8").Value = 'The key to the depths of the mine in Vernis.'
$ws.Range("H8").Value = 'ヴェルニース炭鉱の深層への扉の鍵だ。'

# Row 9
$ws.Range("A9").Value = '501'
$ws.Range("B9").Value = 'Beta 22.57'
$ws.Range("C9").Value = '尼米尔的钥匙'
$ws.Range("D9").Value = 'Key for Nymelle'
$ws.Range("E9").Value = 'ナイミールの鍵'
$ws.Range("F9").Value = '通往尼米尔深层之门的钥匙。'
$ws.Range("G9").Value = 'The key to the depths of Nymelle.'
$ws.Range("H9").Value = 'ナイミール深層への扉の鍵だ。'

# Row 10
$ws.Range("A10").Value = '550'
$ws.Range("B10").Value = 'EA 23.131'
$ws.Range("C10").Value = '宽恕的手环'
$ws.Range("D10").Value = 'Bracelet of Forgiveness'
$ws.Range("E10").Value = '赦しの腕輪'
$ws.Range("F10").Value = '用枯萎的藤蔓编织而成的破旧手环。'
$ws.Range("G10").Value = 'A tattered bracelet woven from withered vines.'
$ws.Range("H10").Value = '枯れた蔦で編まれたボロボロの腕輪だ。'

# Row 11
$ws.Range("A11").Value = '590'
$ws.Range("B11").Value = 'EA 23.50'
$ws.Range("C11").Value = '灯光师执照'
$ws.Range("D11").Value = 'Illumination Engineer License'
$ws.Range("E11").Value = '照明技師のライセンス'
$ws.Range("F11").Value = '灯光师执照。只要带在身上，就能自由改变灯光的颜色。'
$ws.Range("G11").Value = 'A license for Illumination Engineer. Possessing it allows you to freely change colors of lights.'
$ws.Range("H11").Value = '照明技師のライセンスだ。所持することで、照明の色を自由に変えることができるようになる。'

# Row 12
$ws.Range("A12").Value = '600'
$ws.Range("B12").Value = 'Beta 22.57'
$ws.Range("C12").Value = '巢穴探索许可证'
$ws.Range("D12").Value = 'Void Exploration License'
$ws.Range("E12").Value = 'すくつ探索許可証'
$ws.Range("F12").Value = '允许冒险者探索迷之古代遗迹的许可证。'
$ws.Range("G12").Value = 'A license permitting adventurers to explore the mysterious ancient ruins.'
$ws.Range("H12").Value = '冒険者に謎の古代遺跡の探索を許可するライセンスだ。'

# Row 13
$ws.Range("A13").Value = '601'
$ws.Range("B13").Value = 'Beta 22.57'
$ws.Range("C13").Value = '上级冒险者许可证'
$ws.Range("D13").Value = 'Advanced Adventurer License'
$ws.Range("E13").Value = '上級冒険者ライセンス'
$ws.Range("F13").Value = '证明是上级冒险者的许可证。可以获得超越普通冒险者上限的名声，提里斯将生成危险的奈菲亚。'
$ws.Range("G13").Value = 'A license that certifies the holder as an advanced adventurer, allowing acquisition of fame beyond the usual limits for adventurers and results in the generation of dangerous Nefias in Tyris.'
$ws.Range("H13").Value = '上級冒険者であることを証明するライセンスだ。通常の冒険者の上限を超えて名声を獲得できるようになり、ティリスに危険なネフィアが生成されるようになる。'

# Row 14
$ws.Range("A14").Value = '610'
$ws.Range("B14").Value = 'Beta 22.57'
$ws.Range("C14").Value = '幸运硬币'
$ws.Range("D14").Value = 'Lucky Coin'
$ws.Range("E14").Value = 'ラッキーコイン'
$ws.Range("F14").Value = '这是传说中能给持有者带来幸运的硬币。'
$ws.Range("G14").Value = 'A coin said to bring fortune to its possessor.'
$ws.Range("H14").Value = '持つものに幸運をもたらすと噂されるコインだ。'

# Rows 15 and 16 (Bubbling Mysterious Vial 552/553 leftovers) no longer exist; remove them
$ws.Rows("15:16").Delete()
